$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new "PAST_TOBACCO / nonpatient,present" concept-features
# row right after the existing PAST_TOBACCO row (row 7), pushing everything
# below down by one row. Copying row 7 preserves the bold style on column A
# and the column-E "note" style exactly. ---
[void]$ws.Rows.Item(7).EntireRow.Copy()
[void]$ws.Rows.Item(11).EntireRow.Insert()

$ws.Cells.Item(11, 2).Value = "PAST_TOBACCO"
$ws.Cells.Item(11, 7).Value = "Temporality:present"
$ws.Cells.Item(11, 8).Value = "Experiencer:nonpatient"

# --- Step 2: the "COPYALL" data block (originally rows 19-34, now rows
# 20-35 after the insert above) gets reshuffled/rewritten. Clear the whole
# region first and rebuild it from scratch in the target layout. ---
[void]$ws.Range("A20:I35").Clear()

function Set-Row([int]$r, [string]$a, [string]$c, [string]$d, [string]$e)
{
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = "COPYALL"
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    if ($e -ne $null) {
        $ws.Cells.Item($r, 5).Value = $e
    }
}

Set-Row 20 "POSSIBLE_NO_TOBACCO" "Concept" "negated,uncertain" "DocumentAnnotation"
Set-Row 21 "POSSIBLE_TOBACCO"    "Concept" "affirm,uncertain"  "DocumentAnnotation"
Set-Row 22 "HYPOTHETICALCONCEPT" "Concept" "affirm,hypothetical" $null

Set-Row 24 "FAMILY_TOBACCO" "Concept" "family" "DocumentAnnotation"
Set-Row 25 "OTHERS_TOBACCO" "Concept" "others" "DocumentAnnotation"

Set-Row 27 "PAST_TOBACCO" "TOBACCO_CIGARETTES" "affirm,historical,patient" "DocumentAnnotation"
Set-Row 28 "PAST_TOBACCO" "CURRENT_TOBACCO"    "affirm,historical,patient" "DocumentAnnotation"

Set-Row 30 "NO_TOBACCO" "TOBACCO_CIGARETTES" "negated,certain,patient" "DocumentAnnotation"
Set-Row 31 "NO_TOBACCO" "CURRENT_TOBACCO"    "negated,certain,patient" "DocumentAnnotation"
Set-Row 32 "NO_TOBACCO" "PAST_TOBACCO"       "negated"                 "DocumentAnnotation"

Set-Row 33 "CURRENT_TOBACCO" "TOBACCO_CIGARETTES" "affirm,present,certain,patient" "DocumentAnnotation"

# --- Step 3: update the view so the selection matches the saved state
# (top-left back to A1, active cell on A32). ---
[void]$excel.ActiveWindow.ScrollColumn
[void]$ws.Range("A32").Select()
